# Daily report: ensure header rows exist on sheets that have no attendance
# data, and refresh the "Generated" timestamp on the Summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet: bump the "Generated" timestamp (B13)
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B13").Value = "2025-05-19 21:36:44"

# ---------------------------------------------------------------------
# 2. Build a reusable "header" look: bold font (reuse the same bold font
#    already used for the Summary sheet's "Report Sources:" label, A8),
#    a light-gray fill, and centered text - by copying A8's format first
#    and then layering fill + alignment on top.
# ---------------------------------------------------------------------
# ColumnWidth is expressed in characters; Excel stores the serialized
# <col width> in "pixel" units that are offset from the character width
# by the default 5px padding (5/6 of a character at the workbook's
# default font). Subtracting that offset here reproduces the exact
# stored width values (15 and 20) that the export tool wrote.
$narrowWidth = 15 - 0.8333333333333333
$wideWidth   = 20 - 0.8333333333333333
$grayFill    = 14540253   # RGB(221,221,221) = 0xDDDDDD
$centerAlign = -4108      # xlCenter

function Set-HeaderRow {
    param(
        $ws,
        [string[]]$headers
    )

    $n = $headers.Length

    # Column widths: all but the last column are "narrow", last is "wide"
    for ($i = 1; $i -lt $n; $i++) {
        $ws.Columns.Item($i).ColumnWidth = $narrowWidth
    }
    $ws.Columns.Item($n).ColumnWidth = $wideWidth

    # Write header values
    for ($i = 1; $i -le $n; $i++) {
        $ws.Cells.Item(1, $i).Value = $headers[$i - 1]
    }

    # Style: copy the bold "Report Sources:" format from Summary!A8, then
    # apply the header fill + centered alignment on top of it.
    $headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $n))
    $summary.Range("A8").Copy()
    $headerRange.PasteSpecial(-4122)
    $headerRange.Interior.Color = $grayFill
    $headerRange.HorizontalAlignment = $centerAlign
}

# ---------------------------------------------------------------------
# 3. Late Arrivals
# ---------------------------------------------------------------------
$lateArrivals = $wb.Worksheets.Item("Late Arrivals")
Set-HeaderRow $lateArrivals @(
    "Driver", "Asset ID", "Scheduled Start", "Actual Start",
    "Minutes Late", "Job Site", "Division", "Contact Info", "Email"
)

# ---------------------------------------------------------------------
# 4. Early Departures
# ---------------------------------------------------------------------
$earlyDepartures = $wb.Worksheets.Item("Early Departures")
Set-HeaderRow $earlyDepartures @(
    "Driver", "Asset ID", "Scheduled End", "Actual End",
    "Minutes Early", "Job Site", "Division", "Contact Info", "Email"
)

# ---------------------------------------------------------------------
# 5. Not On Job
# ---------------------------------------------------------------------
$notOnJob = $wb.Worksheets.Item("Not On Job")
Set-HeaderRow $notOnJob @(
    "Driver", "Asset ID", "Scheduled Job", "Actual Job",
    "Region", "Division", "Contact Info", "Email"
)

# ---------------------------------------------------------------------
# 6. All Drivers
# ---------------------------------------------------------------------
$allDrivers = $wb.Worksheets.Item("All Drivers")
Set-HeaderRow $allDrivers @(
    "Driver", "Asset ID", "Start Time", "End Time",
    "Total Hours", "Job Site", "Division", "Contact Info", "Email"
)

Write-Host "Header rows applied to all empty attendance sheets."
